$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The previous "Bug" flag text in G1 (written when a timeout/error occurred
# while updating the export) is no longer written now that the timeout
# handling has been fixed. Clear it, but keep the cell's existing style.
$ws.Range("G1").ClearContents()

# Touch the fill so Excel drops the now-unused "applyFill" flag from the
# cell's style (the fill itself stays "none", matching fillId=0).
$ws.Range("G1").Interior.Pattern = 17
$ws.Range("G1").Interior.Pattern = -4142

# The exported value in E2 changed from 999 to 999999.
$ws.Range("E2").Value = 999999

# Leave the cursor on the last-updated cell.
$ws.Range("E2").Select()
